# edit.ps1
# Applies the changes described by the target diff to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rename "Otros gatos" -> "Otros gastos" (typo fix). The shared-string
#    table is managed automatically by the engine: retyping the cell value
#    removes the old unused string and appends the corrected one at the
#    end, which also renumbers every other <v> index referencing the
#    shared-string table exactly like the target workbook.
# ---------------------------------------------------------------------
$ws.Range("G11").Value = "Otros gastos"

# ---------------------------------------------------------------------
# 2) Row 12 (DATOS DE COMPRA figures)
#    - B12 loses its explicit 0 value (now blank)
#    - E12 formula becomes (B12*0.08)*1.21
#    - F12 = 10000, G12 = 500
# ---------------------------------------------------------------------
$ws.Range("B12").ClearContents()
$ws.Range("E12").Formula = "=(B12*0.08)*1.21"
$ws.Range("F12").Value = 10000
$ws.Range("G12").Value = 500

# ---------------------------------------------------------------------
# 3) Row 16 (DATOS DE VENTA figures)
# ---------------------------------------------------------------------
$ws.Range("B16").Value = 10000
$ws.Range("C16").Value = 300
$ws.Range("D16").Value = 600
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 300
$ws.Range("G16").Value = 500

# ---------------------------------------------------------------------
# 4) Row 20 (OTROS DATOS figures)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 0.3
$ws.Range("D20").Value = 25000

# ---------------------------------------------------------------------
# 5) Row 24 (Escenarios de Venta): C24/E24 become blank inputs, D24
#    becomes the average formula and picks up the shaded "computed"
#    style (same style already used by H12/H16/etc, style index 18).
# ---------------------------------------------------------------------
$ws.Range("H12").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Formula = "=(E24+C24)/2"
$ws.Range("C24").ClearContents()
$ws.Range("E24").ClearContents()

# ---------------------------------------------------------------------
# 6) New cell comment on B14 from Alejandro Javier Del Medico Bravo
#    explaining that "DATOS OPERACIÓN" figures are estimates.
# ---------------------------------------------------------------------
$commentText = "Los DATOS OPERACIÓN son aproximados, estimados en base a la mayoria de operaciones.`nEs reponsabilidad de gestor ajustar estos valores de forma precisa para la operación correspondiente."
$ws.Range("B14").AddComment($commentText)

# ---------------------------------------------------------------------
# 7) Page setup: paper size A4->"9" (A4) + portrait orientation.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 8) Selection moves from F24 to B6.
# ---------------------------------------------------------------------
[void]$ws.Range("B6").Select()

Write-Output "edit complete"
